$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.493.46'
$ws.Range('E2').Value = '  -2.62%  '
$ws.Range('D3').Value = '2.477.08'
$ws.Range('E3').Value = '  -1.98%  '
$ws.Range('E4').Value = '  +0.64%  '
$ws.Range('D5').Value = '''314.16'
$ws.Range('E5').Value = '  -0.29%  '
$ws.Range('D6').Value = '''92.95'
$ws.Range('E6').Value = '  -6.10%  '
$ws.Range('E7').Value = '  -3.24%  '
$ws.Range('E8').Value = '  +0.56%  '
$ws.Range('D9').Value = '''0.494'
$ws.Range('E9').Value = '  -4.31%  '
$ws.Range('D10').Value = '''33.11'
$ws.Range('E11').Value = '  -2.68%  '
$ws.Range('E12').Value = '  -0.03%  '
$ws.Range('D13').Value = '2.859.31'
$ws.Range('E13').Value = '  -1.80%  '
$ws.Range('D14').Value = '''6.87'
$ws.Range('E14').Value = '  -4.73%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').Value = '''15.34'
$ws.Range('E15').Value = '  +0.58%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '2.447.18'
$ws.Range('E16').Value = '  -2.76%  '
$ws.Range('E17').Value = '  -3.21%  '
$ws.Range('D18').Value = '41.312.80'
$ws.Range('E18').Value = '  -2.98%  '
$ws.Range('D19').Value = '''6.29'
$ws.Range('E19').Value = '  -4.71%  '
$ws.Range('D20').Value = '0.0₃0922'
$ws.Range('E20').Value = '  -1.79%  '
$ws.Range('D21').Value = '''70.12'
$ws.Range('E21').Value = '  +1.53%  '
$ws.Range('D22').Value = '''11.05'
$ws.Range('E22').Value = '  -9.41%  '
$ws.Range('D23').Value = '''235.09'
$ws.Range('E23').Value = '  -2.86%  '
$ws.Range('D24').Value = '''2.74'
$ws.Range('E24').Value = '  -4.25%  '
$ws.Range('E25').Value = '  +0.03%  '
$ws.Range('D26').Value = '''1.88'
$ws.Range('E26').Value = '  -5.76%  '
$ws.Range('D27').Value = '''24.04'
$ws.Range('E27').Value = '  -6.01%  '
$ws.Range('E28').Value = '  -0.35%  '
$ws.Range('D29').Value = '''9.74'
$ws.Range('E29').Value = '  -2.61%  '
$ws.Range('D30').Value = '''36.47'
$ws.Range('E30').Value = '  -3.36%  '
$ws.Range('D31').Value = '''152.93'
$ws.Range('E32').Value = '  -8.72%  '
$ws.Range('D33').Value = '''2.54'
$ws.Range('E33').Value = '  -3.37%  '
$ws.Range('D34').Value = '''2.54'
$ws.Range('E34').Value = '  -6.29%  '
$ws.Range('D35').Value = '''0.0752'
$ws.Range('E35').Value = '  -3.90%  '
$ws.Range('D36').Value = '''17.68'
$ws.Range('E36').Value = '  +0.86%  '
$ws.Range('D37').Value = '''3.02'
$ws.Range('E37').Value = '  -3.98%  '
$ws.Range('E38').Value = '  -5.90%  '
$ws.Range('E39').Value = '  -3.36%  '
$ws.Range('E40').Value = '  -7.41%  '
$ws.Range('D41').Value = '''4.04'
$ws.Range('E41').Value = '  -4.59%  '
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').Value = '''21.32'
$ws.Range('E42').Value = '  -3.22%  '
$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').Value = '''1.01'
$ws.Range('E43').Value = '  +1.22%  '
$ws.Range('D44').Value = '1.965.43'
$ws.Range('E44').Value = '  -2.52%  '
$ws.Range('D45').Value = '''0.0281'
$ws.Range('E45').Value = '  -4.63%  '
$ws.Range('D46').Value = '''2.96'
$ws.Range('E46').Value = '  -8.02%  '
$ws.Range('D47').Value = '''8.79'
$ws.Range('E47').Value = '  -1.92%  '
$ws.Range('D48').Value = '2.723.09'
$ws.Range('E48').Value = '  -1.45%  '
$ws.Range('D49').Value = '''68.45'
$ws.Range('E49').Value = '  -4.23%  '
$ws.Range('D50').Value = '''95.97'
$ws.Range('E50').Value = '  -4.04%  '
$ws.Range('D51').Value = '''0.176'
$ws.Range('E51').Value = '  -6.44%  '
